$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Model-year refresh: 2020 -> 2021 for a block of existing rows, with a
#    handful of MSRP (column D) bumps mixed in.
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = 2021

$ws.Range("C11").Value = 2021
$ws.Range("D11").Value = 45100

$ws.Range("C12").Value = 2021
$ws.Range("D12").Value = 48900

$ws.Range("C13").Value = 2021

$ws.Range("C14").Value = 2021
$ws.Range("D14").Value = 47010

$ws.Range("C15").Value = 2021
$ws.Range("D15").Value = 50810

$ws.Range("C16").Value = 2021
$ws.Range("D16").Value = 45700

$ws.Range("C55").Value = 2021
$ws.Range("D55").Value = 45070

$ws.Range("C56").Value = 2021
$ws.Range("D56").Value = 46470

$ws.Range("C57").Value = 2021
$ws.Range("D57").Value = 47900

$ws.Range("C58").Value = 2021
$ws.Range("D58").Value = 49300

$ws.Range("C59").Value = 2021
$ws.Range("D59").Value = 48550

$ws.Range("C60").Value = 2021
$ws.Range("D60").Value = 49950

$ws.Range("C61").Value = 2021
$ws.Range("D61").Value = 47720

$ws.Range("C62").Value = 2021
$ws.Range("D62").Value = 51110

$ws.Range("C63").Value = 2021
$ws.Range("D63").Value = 51200

# ---------------------------------------------------------------------------
# 2) Column B is getting much wider to fit the longer trim names below, and
#    the view is being scrolled/zoomed down toward the newly appended rows.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 54.14

$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 80

# ---------------------------------------------------------------------------
# 3) Append the new MSRP rows (81-94), plus a trailing formatted-but-empty
#    row 95. Number formats: column D uses a new no-decimal currency style,
#    column E reuses the existing DPHF currency style.
#
#    The *text* cells (trim names in column B, plus the text-style part
#    codes in column A for the "SE" rows) must be written in this precise
#    order so the rebuilt sharedStrings.xml table lines up exactly with the
#    authored edit - new entries are appended to the shared-string table in
#    first-write order.
# ---------------------------------------------------------------------------
$eFormat = $ws.Range("E2").NumberFormat
$dFormat = "$#,##0_);[Red]($#,##0)"

$ws.Range("B82").Value = "ES 250 LUXURY"
$ws.Range("B83").Value = "ES 250 ULTRA LUXURY"
$ws.Range("B84").Value = "ES 250 F SPORT"
$ws.Range("B85").Value = "ES 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B81").Value = "ES 250"
$ws.Range("B86").Value = "RX 350 F SPORT PERFORMANCE FWD"
$ws.Range("B87").Value = "RX 350 F SPORT PERFORMANCE AWD"
$ws.Range("B88").Value = "RX 450h F-SPORT PERFORMANCE AWD"
$ws.Range("B89").Value = "RX 350L LUXURY FWD"
$ws.Range("B90").Value = "RX 350L LUXURY AWD"
$ws.Range("B91").Value = "RX 450hL LUXURY AWD"
$ws.Range("B92").Value = "RX 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B94").Value = "RX 450h F SPORT AWD BLACK LINE SPECIAL EDITION"
$ws.Range("A92").Value = "9422SE"
$ws.Range("A93").Value = "9426SE"
$ws.Range("A94").Value = "9446SE"
$ws.Range("A85").Value = "9005SE"
$ws.Range("B93").Value = "RX 350 AWD F SPORT BLACK LINE SPECIAL EDITION"

# Remaining (numeric) cells - order doesn't affect the shared-string table.
$ws.Range("A81").Value = 9012
$ws.Range("A82").Value = 9013
$ws.Range("A83").Value = 9014
$ws.Range("A84").Value = 9015
$ws.Range("A86").Value = 9423
$ws.Range("A87").Value = 9427
$ws.Range("A88").Value = 9447
$ws.Range("A89").Value = 9432
$ws.Range("A90").Value = 9436
$ws.Range("A91").Value = 9457

foreach ($row in 81..94) {
    $ws.Range("C$row").Value = 2021
}

$ws.Range("D81").Value = 39900
$ws.Range("D82").Value = 45100
$ws.Range("D83").Value = 48900
$ws.Range("D84").Value = 45700
$ws.Range("D85").Value = 46550
$ws.Range("D86").Value = 50950
$ws.Range("D87").Value = 52350
$ws.Range("D88").Value = 53520
$ws.Range("D89").Value = 53900
$ws.Range("D90").Value = 55300
$ws.Range("D91").Value = 57110
$ws.Range("D92").Value = 49235
$ws.Range("D93").Value = 50635
$ws.Range("D94").Value = 51885

$ws.Range("D81:D94").NumberFormat = $dFormat

foreach ($row in 81..94) {
    $ws.Range("E$row").Value = 1025
    $ws.Range("E$row").NumberFormat = $eFormat
}

# Stray formatted-but-empty cells left behind from the original author's
# copy/paste of row formatting.
$ws.Range("J91").NumberFormat = $eFormat
$ws.Range("K92").NumberFormat = $eFormat
$ws.Range("K93").NumberFormat = $eFormat
$ws.Range("K94").NumberFormat = $eFormat
$ws.Range("K95").NumberFormat = $eFormat

# ---------------------------------------------------------------------------
# 4) Final selection, matching the author's last cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B93").Select()
